$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('A2').Value = 'schubert-winterreise_61'
$ws.Range('B2').Value = 'schubert-winterreise_89'
$ws.Range('C2').Value = 0.1177257525083612
$ws.Range('D2').Value = '[[''C:min/D#'', ''G:min/D'', ''D:7'', ''G:min'']]'
$ws.Range('E2').Value = '[[''C:min/G'', ''G:min'', ''D:7'', ''G:min'']]'
$ws.Range('F2').Value = '[(42.14, 44.36)]'
$ws.Range('G2').Value = '[(2.98, 6.7)]'
$ws.Range('H2').Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'
$ws.Range('I2').Value = 'spotify:track:4lrfYSnZmpXdCWuWqVo8L0'

# Row 3
$ws.Range('A3').Value = 'schubert-winterreise_164'
$ws.Range('B3').Value = 'schubert-winterreise_157'
$ws.Range('C3').Value = 0.2753623188405797
$ws.Range('D3').Value = '[[''A:min'', ''D:min'', ''A:min'', ''E:7'', ''A:min'']]'
$ws.Range('E3').Value = '[[''F:min'', ''A#:min/F'', ''F:min'', ''C:7'', ''F:min'']]'
$ws.Range('F3').Value = '[(14.52, 29.1)]'
$ws.Range('G3').Value = '[(0.3, 4.8)]'
$ws.Range('H3').Value = 'spotify:track:3OD2uwEUQKg0WyW9Lewata'
$ws.Range('I3').Value = 'spotify:track:4lrfYSnZmpXdCWuWqVo8L0'

# Row 4
$ws.Range('A4').Value = 'isophonics_30'
$ws.Range('B4').Value = 'isophonics_107'
$ws.Range('C4').Value = 0.1894736842105263
$ws.Range('D4').Value = '[[''A'', ''A'', ''D'']]'
$ws.Range('E4').Value = '[[''E'', ''E'', ''A'']]'
$ws.Range('F4').Value = '[(14.448702, 17.908476)]'
$ws.Range('G4').Value = '[(122.976598, 130.383764)]'
$ws.Range('H4').Value = ''
$ws.Range('I4').Value = ''

# Row 5
$ws.Range('A5').Value = 'isophonics_297'
$ws.Range('B5').Value = 'isophonics_199'
$ws.Range('C5').Value = 0.06725146198830409
$ws.Range('D5').Value = '[[''C'', ''D'', ''G'']]'
$ws.Range('E5').Value = '[[''G'', ''A'', ''D'']]'
$ws.Range('F5').Value = '[(18.675377, 23.261318)]'
$ws.Range('G5').Value = '[(2.088054, 5.756807)]'
$ws.Range('H5').Value = ''
$ws.Range('I5').Value = ''

# Row 6
$ws.Range('A6').Value = 'schubert-winterreise_69'
$ws.Range('B6').Value = 'schubert-winterreise_6'
$ws.Range('C6').Value = 0.1384615384615385
$ws.Range('D6').Value = '[[''B:min'', ''F#:7/A#'', ''B:min'']]'
$ws.Range('E6').Value = '[[''B:min/F#'', ''F#:7'', ''B:min'']]'
$ws.Range('F6').Value = '[(42.36, 50.36)]'
$ws.Range('G6').Value = '[(79.42, 86.02)]'
$ws.Range('H6').Value = 'spotify:track:1yerCi2iQCVkdHG6rdRn7R'
$ws.Range('I6').Value = 'spotify:track:2g41AZ58LFdQLxmWx82ujI'

# Row 7
$ws.Range('A7').Value = 'isophonics_93'
$ws.Range('B7').Value = 'isophonics_261'
$ws.Range('C7').Value = 0.1154891304347826
$ws.Range('D7').Value = '[[''C'', ''F:maj6'', ''C'']]'
$ws.Range('E7').Value = '[[''G'', ''A:min7'', ''G/3'']]'
$ws.Range('F7').Value = '[(36.36882, 43.729546)]'
$ws.Range('G7').Value = '[(10.107052, 14.565283)]'
$ws.Range('H7').Value = ''
$ws.Range('I7').Value = 'spotify:track:2B4Y9u4ERAFiMo13XPJyGP'

# Row 8
$ws.Range('A8').Value = 'isophonics_31'
$ws.Range('B8').Value = 'isophonics_251'
$ws.Range('C8').Value = 0.3101503759398496
$ws.Range('D8').Value = '[[''E:min'', ''C'', ''G'', ''C'', ''G'']]'
$ws.Range('E8').Value = '[[''B:min'', ''G'', ''D'', ''G'', ''D'']]'
$ws.Range('F8').Value = '[(11.110106, 30.271443)]'
$ws.Range('G8').Value = '[(22.141972, 33.3456)]'
$ws.Range('H8').Value = ''
$ws.Range('I8').Value = ''

# Row 9
$ws.Range('A9').Value = 'isophonics_111'
$ws.Range('B9').Value = 'schubert-winterreise_63'
$ws.Range('C9').Value = 0.1613636363636363
$ws.Range('D9').Value = '[[''C'', ''G:7'', ''C'']]'
$ws.Range('E9').Value = '[[''C/G'', ''G:7'', ''C'']]'
$ws.Range('F9').Value = '[(70.519024, 75.720294)]'
$ws.Range('G9').Value = '[(253.84, 257.6)]'
$ws.Range('H9').Value = ''
$ws.Range('I9').Value = ''

# Row 10
$ws.Range('A10').Value = 'schubert-winterreise_188'
$ws.Range('B10').Value = 'isophonics_212'
$ws.Range('C10').Value = 0.25
$ws.Range('D10').Value = '[[''F:maj'', ''A#:maj'', ''F:maj'', ''A#:maj'', ''F:maj'']]'
$ws.Range('E10').Value = '[[''D'', ''G'', ''D'', ''G'', ''D'']]'
$ws.Range('F10').Value = '[(128.76, 134.6)]'
$ws.Range('G10').Value = '[(46.93228, 57.636679)]'
$ws.Range('H10').Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'
$ws.Range('I10').Value = ''

# Row 11
$ws.Range('A11').Value = 'isophonics_213'
$ws.Range('B11').Value = 'jaah_43'
$ws.Range('C11').Value = 0.1607142857142857
$ws.Range('D11').Value = '[[''E'', ''E:7'', ''A'']]'
$ws.Range('E11').Value = '[[''Eb'', ''Eb:7'', ''Ab'']]'
$ws.Range('F11').Value = '[(16.192174, 27.488682)]'
$ws.Range('G11').Value = '[(45.01, 49.44)]'
$ws.Range('H11').Value = ''
$ws.Range('I11').Value = ''

# Row 12
$ws.Range('A12').Value = 'jaah_49'
$ws.Range('B12').Value = 'schubert-winterreise_156'
$ws.Range('C12').Value = 0.07061307475758523
$ws.Range('D12').Value = '[[''F:7'', ''Bb'', ''Bb'']]'
$ws.Range('E12').Value = '[[''C#:7'', ''F#:maj'', ''F#:maj/A#'']]'
$ws.Range('F12').Value = '[(44.74, 47.14)]'
$ws.Range('G12').Value = '[(21.44, 24.06)]'
$ws.Range('H12').Value = ''
$ws.Range('I12').Value = 'spotify:track:4lrfYSnZmpXdCWuWqVo8L0'

# Row 13
$ws.Range('A13').Value = 'isophonics_0'
$ws.Range('B13').Value = 'jaah_47'
$ws.Range('C13').Value = 0.109375
$ws.Range('D13').Value = '[[''Db:7'', ''Gb'', ''Gb'']]'
$ws.Range('E13').Value = '[[''Db:7'', ''Gb'', ''Gb'']]'
$ws.Range('F13').Value = '[(28.079297, 32.374988)]'
$ws.Range('G13').Value = '[(20.14, 23.58)]'
$ws.Range('H13').Value = ''
$ws.Range('I13').Value = ''

# Row 14
$ws.Range('A14').Value = 'jaah_65'
$ws.Range('B14').Value = 'jaah_62'
$ws.Range('C14').Value = 0.03180619644034278
$ws.Range('D14').Value = '[[''F:7'', ''Bb:min7'', ''Eb:7'']]'
$ws.Range('E14').Value = '[[''D:7'', ''G:min7'', ''C:7'']]'
$ws.Range('F14').Value = '[(7.42, 9.26)]'
$ws.Range('G14').Value = '[(29.42, 30.52)]'
$ws.Range('H14').Value = ''
$ws.Range('I14').Value = ''

# Row 15
$ws.Range('A15').Value = 'isophonics_149'
$ws.Range('B15').Value = 'isophonics_288'
$ws.Range('C15').Value = 0.2094017094017094
$ws.Range('D15').Value = '[[''B'', ''F#'', ''B'', ''E''], [''C#'', ''F#'', ''E'', ''B'']]'
$ws.Range('E15').Value = '[[''E'', ''B'', ''E'', ''A''], [''F#'', ''B'', ''A/9'', ''E/5'']]'
$ws.Range('F15').Value = '[(45.875873, 53.747437), (56.394512, 66.994421)]'
$ws.Range('G15').Value = '[(38.714036, 44.240385), (59.356575, 64.894535)]'
$ws.Range('H15').Value = ''
$ws.Range('I15').Value = ''

# Row 16
$ws.Range('A16').Value = 'isophonics_61'
$ws.Range('B16').Value = 'schubert-winterreise_14'
$ws.Range('C16').Value = 0.162280701754386
$ws.Range('D16').Value = '[[''F#:min'', ''C#'', ''F#:min'']]'
$ws.Range('E16').Value = '[[''D:min'', ''A:maj'', ''D:min'']]'
$ws.Range('F16').Value = '[(3.744784, 6.774988)]'
$ws.Range('G16').Value = '[(10.56, 17.34)]'
$ws.Range('H16').Value = ''
$ws.Range('I16').Value = ''

# Row 17
$ws.Range('A17').Value = 'schubert-winterreise_88'
$ws.Range('B17').Value = 'schubert-winterreise_55'
$ws.Range('C17').Value = 0.2657342657342657
$ws.Range('D17').Value = '[[''A:maj/E'', ''E:7'', ''A:maj'', ''E:7'', ''A:maj'']]'
$ws.Range('E17').Value = '[[''G:maj'', ''D:7/C'', ''G:maj/B'', ''D:7/C'', ''G:maj/B'']]'
$ws.Range('F17').Value = '[(16.04, 21.0)]'
$ws.Range('G17').Value = '[(68.64, 83.84)]'
$ws.Range('H17').Value = 'spotify:track:0XfunCHFEeQnzm4NaY8rJr'
$ws.Range('I17').Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'
